$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2022" column (S), mirroring the formatting of
# the existing "2021" column (R) for every row of the table (rows 3-34).
$ws.Range("R3").Copy($ws.Range("S3"))
$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("R5").Copy($ws.Range("S5"))
$ws.Range("R6").Copy($ws.Range("S6"))
$ws.Range("R7").Copy($ws.Range("S7"))
$ws.Range("R8").Copy($ws.Range("S8"))
$ws.Range("R9").Copy($ws.Range("S9"))
$ws.Range("R10").Copy($ws.Range("S10"))
$ws.Range("R11").Copy($ws.Range("S11"))
$ws.Range("R12").Copy($ws.Range("S12"))
$ws.Range("R13").Copy($ws.Range("S13"))
$ws.Range("R14").Copy($ws.Range("S14"))
$ws.Range("R15").Copy($ws.Range("S15"))
$ws.Range("R16").Copy($ws.Range("S16"))
$ws.Range("R17").Copy($ws.Range("S17"))
$ws.Range("R18").Copy($ws.Range("S18"))
$ws.Range("R19").Copy($ws.Range("S19"))
$ws.Range("R20").Copy($ws.Range("S20"))
$ws.Range("R21").Copy($ws.Range("S21"))
$ws.Range("R22").Copy($ws.Range("S22"))
$ws.Range("R23").Copy($ws.Range("S23"))
$ws.Range("R24").Copy($ws.Range("S24"))
$ws.Range("R25").Copy($ws.Range("S25"))
$ws.Range("R26").Copy($ws.Range("S26"))
$ws.Range("R27").Copy($ws.Range("S27"))
$ws.Range("R28").Copy($ws.Range("S28"))
$ws.Range("R29").Copy($ws.Range("S29"))
$ws.Range("R30").Copy($ws.Range("S30"))
$ws.Range("R31").Copy($ws.Range("S31"))
$ws.Range("R32").Copy($ws.Range("S32"))
$ws.Range("R33").Copy($ws.Range("S33"))
$ws.Range("R34").Copy($ws.Range("S34"))

# Fill in the 2022 figures (row 4 is the year header, rows 5-34 the data;
# row 3 is the blank separator row above the header and keeps no value).
$ws.Range("S4").Value = 2022

$ws.Range("S5").Value = 135
$ws.Range("S6").Value = 99
$ws.Range("S7").Value = 36

$ws.Range("S8").Value = 97
$ws.Range("S9").Value = 80
$ws.Range("S10").Value = 17

$ws.Range("S11").Value = 17
$ws.Range("S12").Value = 11
$ws.Range("S13").Value = 6

$ws.Range("S14").Value = 5
$ws.Range("S15").Value = 3
$ws.Range("S16").Value = 2

$ws.Range("S17").Value = "-"
$ws.Range("S18").Value = "-"
$ws.Range("S19").Value = "-"

$ws.Range("S20").Value = 6
$ws.Range("S21").Value = 1
$ws.Range("S22").Value = 5

$ws.Range("S23").Value = "-"
$ws.Range("S24").Value = "-"
$ws.Range("S25").Value = "-"

$ws.Range("S26").Value = 10
$ws.Range("S27").Value = 4
$ws.Range("S28").Value = 6

$ws.Range("S29").Value = "-"
$ws.Range("S30").Value = "-"
$ws.Range("S31").Value = "-"

$ws.Range("S32").Value = "-"
$ws.Range("S33").Value = "-"
$ws.Range("S34").Value = "-"

# Selected cell moves on to the next empty column, matching the authored file.
$ws.Range("T4").Select() | Out-Null
